$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add the new "topic_link" column (H) ---

# Header cell H1: copy the header formatting from the existing bold header
# cell G1 ("webinar_topic") so H1 matches the other header cells, then set
# its text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "topic_link"

# Data rows: one relative link per webinar topic row (rows 2-4 have no
# webinar topic link, so they are left blank as in the source data).
# Values are entered in the same order the author typed them in (the four
# Shiny sessions first, then the rest top-to-bottom) so the shared-string
# table is built up in the same sequence as the original edit.
$ws.Range("H5").Value  = "topics/shiny_1.html"
$ws.Range("H7").Value  = "topics/shiny_2.html"
$ws.Range("H11").Value = "topics/shiny_4.html"
$ws.Range("H9").Value  = "topics/shiny_3.html"
$ws.Range("H6").Value  = "topics/tidyverse.html"
$ws.Range("H8").Value  = "topics/geospatial.html"
$ws.Range("H10").Value = "topics/census.html"
$ws.Range("H12").Value = "topics/machine_learning.html"
$ws.Range("H13").Value = "topics/rmarkdown.html"
$ws.Range("H14").Value = "topics/misc_questions.html"

# Match the column width used for the new column in the target workbook.
$ws.Columns.Item(8).ColumnWidth = 24.1

# Match the selection left behind in the saved workbook.
$ws.Range("H12").Select()

Write-Host "Applied topic_link column changes"
